# Applies the weekly crime-data refresh described in the commit
# "New crime data collected" to the CompStat_1 worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and report week dates) ---
$ws.Range("A8").Value = "Volume 30   Number  37"
$ws.Range("C9").Value = "Report Covering the Week  9/11/2023  Through  9/17/2023"

# --- Cells that change data type (count <-> "N/A" / "***.*") ---
# Copy formatting from an untouched same-type cell first so the
# resulting style index matches a genuine number/text style,
# then set the real value.
$ws.Range("I14").Copy($ws.Range("C18"))
$ws.Range("C18").Value = 2
$ws.Range("J14").Copy($ws.Range("C22"))
$ws.Range("C22").Value = 3
$ws.Range("C14").Copy($ws.Range("D28"))
$ws.Range("D14").Copy($ws.Range("D29"))
$ws.Range("E14").Copy($ws.Range("E28"))
$ws.Range("H14").Copy($ws.Range("E29"))

# --- Remaining numeric value updates (rows 14-29) ---
$ws.Range("L14").Value = -57.142857142857
$ws.Range("N14").Value = -88.888888888888
$ws.Range("D15").Value = 1
$ws.Range("F15").Value = 2
$ws.Range("H15").Value = -50
$ws.Range("J15").Value = 25
$ws.Range("K15").Value = 12
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 9
$ws.Range("E16").Value = -66.666666666666
$ws.Range("F16").Value = 13
$ws.Range("G16").Value = 23
$ws.Range("H16").Value = -43.478260869565
$ws.Range("I16").Value = 188
$ws.Range("J16").Value = 254
$ws.Range("K16").Value = -25.984251968503
$ws.Range("L16").Value = 27.027027027027
$ws.Range("M16").Value = -40.506329113924
$ws.Range("N16").Value = -88.359133126935
$ws.Range("C17").Value = 15
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = 87.5
$ws.Range("F17").Value = 57
$ws.Range("G17").Value = 56
$ws.Range("H17").Value = 1.785714285714
$ws.Range("I17").Value = 520
$ws.Range("J17").Value = 519
$ws.Range("K17").Value = 0.192678227360
$ws.Range("L17").Value = 18.721461187214
$ws.Range("M17").Value = 15.555555555555
$ws.Range("N17").Value = -37.799043062201
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 9
$ws.Range("G18").Value = 10
$ws.Range("H18").Value = -10
$ws.Range("I18").Value = 87
$ws.Range("J18").Value = 138
$ws.Range("K18").Value = -36.956521739130
$ws.Range("L18").Value = -8.421052631578
$ws.Range("M18").Value = -57.766990291262
$ws.Range("N18").Value = -84.974093264248
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = -14.285714285714
$ws.Range("F19").Value = 25
$ws.Range("G19").Value = 27
$ws.Range("H19").Value = -7.407407407407
$ws.Range("I19").Value = 243
$ws.Range("J19").Value = 283
$ws.Range("K19").Value = -14.134275618374
$ws.Range("L19").Value = 19.704433497536
$ws.Range("M19").Value = -22.611464968152
$ws.Range("N19").Value = -64
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 9
$ws.Range("H20").Value = -43.75
$ws.Range("I20").Value = 111
$ws.Range("J20").Value = 137
$ws.Range("K20").Value = -18.978102189781
$ws.Range("L20").Value = 8.823529411764
$ws.Range("M20").Value = -0.892857142857
$ws.Range("N20").Value = -75.442477876106
$ws.Range("C21").Value = 29
$ws.Range("D21").Value = 30
$ws.Range("E21").Value = -3.333333333333
$ws.Range("F21").Value = 115
$ws.Range("G21").Value = 136
$ws.Range("H21").Value = -15.441176470588
$ws.Range("I21").Value = 1183
$ws.Range("J21").Value = 1369
$ws.Range("K21").Value = -13.586559532505
$ws.Range("L21").Value = 14.965986394557
$ws.Range("M21").Value = -17.675713291579
$ws.Range("N21").Value = -72.321010762751
$ws.Range("F22").Value = 6
$ws.Range("H22").Value = 500
$ws.Range("I22").Value = 35
$ws.Range("K22").Value = -5.405405405405
$ws.Range("L22").Value = 34.615384615384
$ws.Range("M22").Value = -7.894736842105
$ws.Range("C23").Value = 7
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = 250
$ws.Range("F23").Value = 25
$ws.Range("G23").Value = 29
$ws.Range("H23").Value = -13.793103448275
$ws.Range("I23").Value = 254
$ws.Range("J23").Value = 277
$ws.Range("K23").Value = -8.303249097472
$ws.Range("L23").Value = 12.389380530973
$ws.Range("M23").Value = 42.696629213483
$ws.Range("C24").Value = 14
$ws.Range("D24").Value = 21
$ws.Range("E24").Value = -33.333333333333
$ws.Range("F24").Value = 63
$ws.Range("G24").Value = 76
$ws.Range("H24").Value = -17.105263157894
$ws.Range("I24").Value = 797
$ws.Range("J24").Value = 817
$ws.Range("K24").Value = -2.447980416156
$ws.Range("L24").Value = 46.507352941176
$ws.Range("M24").Value = 14.020028612303
$ws.Range("C25").Value = 26
$ws.Range("D25").Value = 23
$ws.Range("E25").Value = 13.043478260869
$ws.Range("F25").Value = 97
$ws.Range("G25").Value = 71
$ws.Range("H25").Value = 36.619718309859
$ws.Range("I25").Value = 715
$ws.Range("J25").Value = 689
$ws.Range("K25").Value = 3.773584905660
$ws.Range("L25").Value = 24.564459930313
$ws.Range("M25").Value = -26.288659793814
$ws.Range("D26").Value = 1
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 4
$ws.Range("G26").Value = 4
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 35
$ws.Range("J26").Value = 39
$ws.Range("K26").Value = -10.256410256410
$ws.Range("L26").Value = -28.571428571428
$ws.Range("C27").Value = 2
$ws.Range("D27").Value = 3
$ws.Range("E27").Value = -33.333333333333
$ws.Range("F27").Value = 7
$ws.Range("G27").Value = 7
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 64
$ws.Range("J27").Value = 63
$ws.Range("K27").Value = 1.587301587301
$ws.Range("L27").Value = 1.587301587301
$ws.Range("G28").Value = 5
$ws.Range("L28").Value = -47.540983606557
$ws.Range("N28").Value = -86.666666666666
$ws.Range("G29").Value = 4
$ws.Range("L29").Value = -40
$ws.Range("N29").Value = -86.486486486486
